$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DateOfBirth column (F) had a typo in row 2 ("12-03-2o11") and the
# two test rows' dates were effectively swapped versus row 3 ("01-05-2014").
# Fix the typo and put the corrected date into F2, keep F3 as-is; both
# stay text (leading apostrophe = quotePrefix) so Excel doesn't silently
# turn them into serial date numbers.
$ws.Range("F2").Value = "'12-03-2011"
$ws.Range("F3").Value = "'01-05-2014"

# Give the corrected date cell an explicit date number format.
$ws.Range("F2").NumberFormat = "m/d/yyyy"

# Move the active selection to F8.
[void]$ws.Range("F8").Select()
